$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2831691422999417
$ws.Range("C2").Value = 0.994552455368806
$ws.Range("D2").Value = 0.4309702347941847
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5))])"
$ws.Range("G2").Value = 0.1237476138499915
$ws.Range("H2").Value = 0.992
